$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-6 (columns D, M, N, O, P, R, S), representing a
# re-shuffle of the market data rows (weekly reporting reorder).
# Row mapping (new row <- old row): 2<-5, 3<-4, 4<-6, 5<-2, 6<-3

$data = @{
    2 = @{ D = 44320; M = 50;  N = 18000; O = 20000; P = 18800; R = "Provincia de Limarí";   S = 1044 }
    3 = @{ D = 44719; M = 50;  N = 20000; O = 21000; P = 20400; R = "Provincia de Limarí";   S = 1133 }
    4 = @{ D = 44362; M = 100; N = 19000; O = 20000; P = 19500; R = "Provincia de Curicó";   S = 1083 }
    5 = @{ D = 45084; M = 100; N = 17000; O = 18000; P = 17500; R = "Región de O'Higgins";   S = 972  }
    6 = @{ D = 45106; M = 50;  N = 10000; O = 10000; P = 10000; R = "Región de O'Higgins";   S = 556  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
